$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 462, shifting existing rows 462:488 down to 463:489
$ws.Rows(462).Insert()

# Populate the newly inserted row 462 with the new record's data
$ws.Cells.Item(462, 1).Value = 5
$ws.Cells.Item(462, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(462, 3).Value = "Maule"
$ws.Cells.Item(462, 4).Value = 44826
$ws.Cells.Item(462, 5).Value = 7
$ws.Cells.Item(462, 6).Value = 100112043
$ws.Cells.Item(462, 7).Value = "Pepino ensalada"
$ws.Cells.Item(462, 8).Value = "Sin especificar"
$ws.Cells.Item(462, 9).Value = "Primera"
$ws.Cells.Item(462, 10).Value = 400
$ws.Cells.Item(462, 11).Value = 18000
$ws.Cells.Item(462, 12).Value = 18000
$ws.Cells.Item(462, 13).Value = 18000
$ws.Cells.Item(462, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(462, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(462, 16).Value = 300
$ws.Cells.Item(462, 17).Value = 60
$ws.Cells.Item(462, 18).Value = "Hortaliza"

# Ensure the date cell keeps the existing date-style format used by column D
$ws.Cells.Item(462, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
